# Reorganize regression table:
#  - rename the "rincvar" columns to "inciqr"
#  - move the "educ_gr=low educ" row up, right after "HHinc_gr=low inc",
#    shifting the age-group rows down by one row-pair
#  - refresh the numeric results that came out of the re-run regressions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row: rincvar I/II/III -> inciqr I/II/III ---
$ws.Range("E1").Value = "inciqr I"
$ws.Range("F1").Value = "inciqr II"
$ws.Range("G1").Value = "inciqr III"

# --- rmse row (row 2) + its se row (row 3) ---
$ws.Range("B2").Value = "4.02***"
$ws.Range("E2").Value = "1.97***"
$ws.Range("F2").Value = "1.99***"
$ws.Range("G2").Value = "1.05***"

$ws.Range("E3").Value = "(0.13)"
$ws.Range("F3").Value = "(0.13)"
$ws.Range("G3").Value = "(0.28)"

# --- HHinc_gr=low inc row (row 4) + its se row (row 5) ---
$ws.Range("B4").Value = "0.15***"
$ws.Range("E4").Value = "0.20***"
$ws.Range("F4").Value = "0.20***"
$ws.Range("G4").Value = "0.19***"

# (row 5 standard errors are unchanged)

# --- educ_gr=low educ row, now placed at rows 6-7 ---
# F6 ("0.01") looks like a plain number, so Excel would silently store it
# as a numeric cell instead of text unless the cell is pre-formatted as
# Text; force Text, write the value, then drop back to the default style
# so no stray formatting is left behind on the cell.
$ws.Range("A6").Value = "educ_gr=low educ"
$ws.Range("C6").Value = "-0.16***"
$ws.Range("D6").Value = "-0.12***"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "0.01"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "0.03***"

$ws.Range("C7").Value = "(0.02)"
$ws.Range("D7").Value = "(0.02)"
$ws.Range("F7").Value = "(0.01)"
$ws.Range("G7").Value = "(0.01)"

# --- age_gr=30-39 row, shifted down to rows 8-9 ---
$ws.Range("A8").Value = "age_gr=30-39"
$ws.Range("D8").Value = "-0.32***"
$ws.Range("G8").Value = "-0.16***"

$ws.Range("G9").Value = "(0.01)"

# --- age_gr=40-48 row, shifted down to rows 10-11 ---
$ws.Range("A10").Value = "age_gr=40-48"
$ws.Range("D10").Value = "-0.48***"
$ws.Range("G10").Value = "-0.24***"

$ws.Range("G11").Value = "(0.01)"

# --- age_gr=49-57 row, shifted down to rows 12-13 ---
$ws.Range("A12").Value = "age_gr=49-57"
$ws.Range("D12").Value = "-0.58***"
$ws.Range("G12").Value = "-0.29***"

$ws.Range("D13").Value = "(0.03)"
$ws.Range("G13").Value = "(0.02)"

# --- age_gr=>57 row, shifted down to rows 14-15 (replaces the old educ row) ---
$ws.Range("A14").Value = "age_gr=>57"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = "-0.45***"
$ws.Range("E14").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("G14").Value = "-0.23***"

$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("D15").Value = "(0.04)"
$ws.Range("E15").ClearContents()
$ws.Range("F15").ClearContents()

# --- N row (row 16): the inciqr sample grew ---
# Same plain-number-looking-text situation as F6 above for E16:G17.
$ws.Range("E16:G17").NumberFormat = "@"

$ws.Range("E16").Value = "44874"
$ws.Range("F16").Value = "44874"
$ws.Range("G16").Value = "44874"

# --- R2 row (row 17) ---
$ws.Range("E17").Value = "0.02"
$ws.Range("F17").Value = "0.02"
$ws.Range("G17").Value = "0.03"

$ws.Range("E16:G17").Style = "Normal"
